# Refresh market-price-derived leve-profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H:N) across all crafter sheets, pulling in
# updated Universalis price data for the affected item rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce / Rubber
$ws.Cells.Item(11, 8).Value = 106.77778
$ws.Cells.Item(11, 9).Value = 106.77778
$ws.Cells.Item(11, 11).Value = 106.77778
$ws.Cells.Item(11, 13).Value = 33.22221999999999

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation / Holy Water
$ws.Cells.Item(70, 8).Value = 5397.421
$ws.Cells.Item(70, 9).Value = 2825
$ws.Cells.Item(70, 10).Value = 6083.4
$ws.Cells.Item(70, 11).Value = 8475
$ws.Cells.Item(70, 12).Value = 18250.2
$ws.Cells.Item(70, 13).Value = -8205
$ws.Cells.Item(70, 14).Value = -18790.2

$ws = $wb.Worksheets.Item("ALC")
# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Cells.Item(73, 8).Value = 5397.421
$ws.Cells.Item(73, 9).Value = 2825
$ws.Cells.Item(73, 10).Value = 6083.4
$ws.Cells.Item(73, 11).Value = 8475
$ws.Cells.Item(73, 12).Value = 18250.2
$ws.Cells.Item(73, 13).Value = -7539
$ws.Cells.Item(73, 14).Value = -20122.2

$ws = $wb.Worksheets.Item("ALC")
# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Cells.Item(80, 8).Value = 794.8946999999999
$ws.Cells.Item(80, 9).Value = 541.8333
$ws.Cells.Item(80, 10).Value = 1228.7142
$ws.Cells.Item(80, 11).Value = 1625.4999
$ws.Cells.Item(80, 12).Value = 3686.1426
$ws.Cells.Item(80, 13).Value = -627.4999
$ws.Cells.Item(80, 14).Value = -5682.142599999999

$ws = $wb.Worksheets.Item("ALC")
# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Cells.Item(83, 8).Value = 794.8946999999999
$ws.Cells.Item(83, 9).Value = 541.8333
$ws.Cells.Item(83, 10).Value = 1228.7142
$ws.Cells.Item(83, 11).Value = 4876.4997
$ws.Cells.Item(83, 12).Value = 11058.4278
$ws.Cells.Item(83, 13).Value = 115.5002999999997
$ws.Cells.Item(83, 14).Value = -21042.4278

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 1826.4359
$ws.Cells.Item(132, 9).Value = 1392.3429
$ws.Cells.Item(132, 11).Value = 4177.028700000001
$ws.Cells.Item(132, 13).Value = -1647.028700000001

$ws = $wb.Worksheets.Item("ALC")
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 2070.973
$ws.Cells.Item(137, 9).Value = 1891.32
$ws.Cells.Item(137, 10).Value = 2445.25
$ws.Cells.Item(137, 11).Value = 5673.96
$ws.Cells.Item(137, 12).Value = 7335.75
$ws.Cells.Item(137, 13).Value = -3123.96
$ws.Cells.Item(137, 14).Value = -12435.75

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Cells.Item(5, 8).Value = 5922.1113
$ws.Cells.Item(5, 9).Value = 6624.875
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 6624.875
$ws.Cells.Item(5, 12).Value = 300
$ws.Cells.Item(5, 13).Value = -6512.875
$ws.Cells.Item(5, 14).Value = -524

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value = 2629.1765
$ws.Cells.Item(32, 9).Value = 895.45
$ws.Cells.Item(32, 11).Value = 895.45
$ws.Cells.Item(32, 13).Value = -608.45

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 3344.8386
$ws.Cells.Item(61, 9).Value = 2523.2886
$ws.Cells.Item(61, 10).Value = 7616.9
$ws.Cells.Item(61, 11).Value = 2523.2886
$ws.Cells.Item(61, 12).Value = 7616.9
$ws.Cells.Item(61, 13).Value = -2311.2886
$ws.Cells.Item(61, 14).Value = -8040.9

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 1953.0358
$ws.Cells.Item(74, 9).Value = 1941.9615
$ws.Cells.Item(74, 10).Value = 2097
$ws.Cells.Item(74, 11).Value = 1941.9615
$ws.Cells.Item(74, 12).Value = 2097
$ws.Cells.Item(74, 13).Value = -1067.9615
$ws.Cells.Item(74, 14).Value = -3845

$ws = $wb.Worksheets.Item("ARM")
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 1953.0358
$ws.Cells.Item(77, 9).Value = 1941.9615
$ws.Cells.Item(77, 10).Value = 2097
$ws.Cells.Item(77, 11).Value = 9709.807499999999
$ws.Cells.Item(77, 12).Value = 10485
$ws.Cells.Item(77, 13).Value = -5341.807499999999
$ws.Cells.Item(77, 14).Value = -19221

$ws = $wb.Worksheets.Item("ARM")
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 5569
$ws.Cells.Item(102, 9).Value = 1371.091
$ws.Cells.Item(102, 10).Value = 9417.083000000001
$ws.Cells.Item(102, 11).Value = 1371.091
$ws.Cells.Item(102, 12).Value = 9417.083000000001
$ws.Cells.Item(102, 13).Value = 250.9090000000001
$ws.Cells.Item(102, 14).Value = -12661.083

$ws = $wb.Worksheets.Item("ARM")
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 3344.8386
$ws.Cells.Item(136, 9).Value = 2523.2886
$ws.Cells.Item(136, 10).Value = 7616.9
$ws.Cells.Item(136, 11).Value = 7569.8658
$ws.Cells.Item(136, 12).Value = 22850.7
$ws.Cells.Item(136, 13).Value = -5019.8658
$ws.Cells.Item(136, 14).Value = -27950.7

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Cells.Item(4, 8).Value = 5922.1113
$ws.Cells.Item(4, 9).Value = 6624.875
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 6624.875
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = -6509.875
$ws.Cells.Item(4, 14).Value = -530

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 2101.7
$ws.Cells.Item(99, 9).Value = 985.8
$ws.Cells.Item(99, 10).Value = 3217.6
$ws.Cells.Item(99, 11).Value = 985.8
$ws.Cells.Item(99, 12).Value = 3217.6
$ws.Cells.Item(99, 13).Value = 512.2
$ws.Cells.Item(99, 14).Value = -6213.6

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 1533.9
$ws.Cells.Item(134, 9).Value = 1367.16
$ws.Cells.Item(134, 10).Value = 2367.6
$ws.Cells.Item(134, 11).Value = 4101.48
$ws.Cells.Item(134, 12).Value = 7102.799999999999
$ws.Cells.Item(134, 13).Value = -1566.48
$ws.Cells.Item(134, 14).Value = -12172.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 3819.0527
$ws.Cells.Item(31, 9).Value = 1871.65
$ws.Cells.Item(31, 11).Value = 1871.65
$ws.Cells.Item(31, 13).Value = -1576.65

$ws = $wb.Worksheets.Item("CRP")
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 3819.0527
$ws.Cells.Item(34, 9).Value = 1871.65
$ws.Cells.Item(34, 11).Value = 1871.65
$ws.Cells.Item(34, 13).Value = -1669.65

$ws = $wb.Worksheets.Item("CRP")
# Row 45: A Tree Grew in Gridania / Pastoral Oak Cane
$ws.Cells.Item(45, 8).Value = 5537
$ws.Cells.Item(45, 9).Value = 10000
$ws.Cells.Item(45, 10).Value = 1074
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 12).Value = 1074
$ws.Cells.Item(45, 13).Value = -9407
$ws.Cells.Item(45, 14).Value = -2260

$ws = $wb.Worksheets.Item("CRP")
# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Cells.Item(80, 8).Value = 64999.5
$ws.Cells.Item(80, 10).Value = 64999.5
$ws.Cells.Item(80, 12).Value = 64999.5
$ws.Cells.Item(80, 14).Value = -67245.5

$ws = $wb.Worksheets.Item("CRP")
# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Cells.Item(83, 8).Value = 64999.5
$ws.Cells.Item(83, 10).Value = 64999.5
$ws.Cells.Item(83, 12).Value = 194998.5
$ws.Cells.Item(83, 14).Value = -206230.5

$ws = $wb.Worksheets.Item("CUL")
# Row 32: Convalescence Precedes Essence / Ginger Cookie
$ws.Cells.Item(32, 8).Value = 472
$ws.Cells.Item(32, 9).Value = 1399
$ws.Cells.Item(32, 10).Value = 356.125
$ws.Cells.Item(32, 11).Value = 4197
$ws.Cells.Item(32, 12).Value = 1068.375
$ws.Cells.Item(32, 13).Value = -3914
$ws.Cells.Item(32, 14).Value = -1634.375

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch / Chamomile Tea
$ws.Cells.Item(34, 8).Value = 2301.7222
$ws.Cells.Item(34, 10).Value = 2453.2
$ws.Cells.Item(34, 12).Value = 7359.599999999999
$ws.Cells.Item(34, 14).Value = -7527.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 39: Bloody Good Tart, This / Blood Currant Tart
$ws.Cells.Item(39, 8).Value = 4966.3335
$ws.Cells.Item(39, 10).Value = 4966.3335
$ws.Cells.Item(39, 12).Value = 14899.0005
$ws.Cells.Item(39, 14).Value = -15487.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 55: Pagan Pastries / Pastry Fish
$ws.Cells.Item(55, 8).Value = 2084404.5
$ws.Cells.Item(55, 9).Value = 625744.25
$ws.Cells.Item(55, 10).Value = 5001725
$ws.Cells.Item(55, 11).Value = 1877232.75
$ws.Cells.Item(55, 12).Value = 15005175
$ws.Cells.Item(55, 13).Value = -1877055.75
$ws.Cells.Item(55, 14).Value = -15005529

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Cells.Item(102, 8).Value = 1273.3914
$ws.Cells.Item(102, 9).Value = 656
$ws.Cells.Item(102, 10).Value = 3496
$ws.Cells.Item(102, 11).Value = 656
$ws.Cells.Item(102, 12).Value = 3496
$ws.Cells.Item(102, 13).Value = 966
$ws.Cells.Item(102, 14).Value = -6740

$ws = $wb.Worksheets.Item("LTW")
# Row 80: Don't Sweat the Small Fry / Dragonskin Wristbands
$ws.Cells.Item(80, 8).Value = 68900
$ws.Cells.Item(80, 10).Value = 68900
$ws.Cells.Item(80, 12).Value = 68900
$ws.Cells.Item(80, 14).Value = -71146

$ws = $wb.Worksheets.Item("LTW")
# Row 83: It's All in the Wrists (L) / Dragonskin Wristbands
$ws.Cells.Item(83, 8).Value = 68900
$ws.Cells.Item(83, 10).Value = 68900
$ws.Cells.Item(83, 12).Value = 206700
$ws.Cells.Item(83, 14).Value = -217932

$ws = $wb.Worksheets.Item("LTW")
# Row 96: Off the Cuff / Gyuki Leather Wristband
$ws.Cells.Item(96, 8).Value = 42500
$ws.Cells.Item(96, 10).Value = 42500
$ws.Cells.Item(96, 12).Value = 42500
$ws.Cells.Item(96, 14).Value = -47992

$ws = $wb.Worksheets.Item("LTW")
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Cells.Item(136, 8).Value = 12334.529
$ws.Cells.Item(136, 9).Value = 1954.2
$ws.Cells.Item(136, 10).Value = 35041.5
$ws.Cells.Item(136, 11).Value = 5862.6
$ws.Cells.Item(136, 12).Value = 105124.5
$ws.Cells.Item(136, 13).Value = -3312.6
$ws.Cells.Item(136, 14).Value = -110224.5

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 2878.1482
$ws.Cells.Item(136, 9).Value = 1718.8572
$ws.Cells.Item(136, 11).Value = 5156.571599999999
$ws.Cells.Item(136, 13).Value = -2606.571599999999
